# Weekly crypto price/volume refresh (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as "29.024.82" or "0.4590" that look
# numeric to Excel's auto-detection and would otherwise be coerced into
# floating point numbers (losing the literal text / trailing zeros).
# Temporarily force the column to Text format while writing the new values,
# then clear the format again so the cells end up back at the default style
# (matching how the sheet was originally authored) while keeping the text value.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.991.03"
$ws.Range("D3").Value = "1.924.44"
$ws.Range("D4").Value = "1.006"
$ws.Range("D5").Value = "325.17"
$ws.Range("D7").Value = "0.4583"
$ws.Range("D8").Value = "0.3816"
$ws.Range("D9").Value = "0.07753"
$ws.Range("D10").Value = "0.9792"
$ws.Range("D12").Value = "1.967.49"
$ws.Range("D13").Value = "5.711"
$ws.Range("D14").Value = "6.984"
$ws.Range("D15").Value = "0.07000"
$ws.Range("D16").Value = "84.85"
$ws.Range("D17").Value = "1.005"
$ws.Range("D18").Value = "0.000009503"
$ws.Range("D19").Value = "16.72"
$ws.Range("D21").Value = "29.031.88"
$ws.Range("D22").Value = "5.353"
$ws.Range("D23").Value = "11.07"
$ws.Range("D24").Value = "2.208.49"
$ws.Range("D25").Value = "2.056"
$ws.Range("D26").Value = "158.04"
$ws.Range("D27").Value = "19.01"
$ws.Range("D28").Value = "5.639"
$ws.Range("D29").Value = "117.55"
$ws.Range("D30").Value = "1.838"
$ws.Range("D31").Value = "0.09324"
$ws.Range("D32").Value = "0.8648"
$ws.Range("D33").Value = "5.110"
$ws.Range("D34").Value = "1.247"
$ws.Range("D36").Value = "0.05707"
$ws.Range("D37").Value = "1.151"
$ws.Range("D39").Value = "0.02053"
$ws.Range("D40").Value = "3.092"
$ws.Range("D41").Value = "7.476"
$ws.Range("D42").Value = "0.5516"
$ws.Range("D44").Value = "9.376"
$ws.Range("D45").Value = "0.000002845"
$ws.Range("D46").Value = "2.178"
$ws.Range("D47").Value = "0.5195"
$ws.Range("D48").Value = "0.06944"
$ws.Range("D49").Value = "11.25"
$ws.Range("D50").Value = "111.08"
$ws.Range("D51").Value = "1.764"

$dRange.ClearFormats()

# Columns B, C (coin name / link) and E (already-textual percentage strings)
# are safe to assign directly.
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("E12").Value = "  +2.73%  "
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("E24").Value = "  +2.56%  "
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("E40").Value = "  +13.25%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("E45").Value = "  +10.00%  "
$ws.Range("E46").Value = "  +4.55%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("E51").Value = "  -0.69%  "

